# "Generate Report for Archive"
# The status "Ready for handoff" moves on to "In Translation" for the
# zh-cn and de-de locales. That status string shows up in three places:
#   - Overview sheet: columns E (zh-cn) and F (de-de) of the data row
#   - zh-cn sheet: the "Status" column (C) of the data row
#   - de-de sheet: the "Status" column (C) of the data row

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The shorter status text lets the Status-related columns narrow a bit;
# mirror that column-width shrink on the affected sheets/columns.
$overview.Columns("E:F").ColumnWidth = 12.576851254417766
$zhcn.Columns("C:C").ColumnWidth = 12.576851254417766
$dede.Columns("C:C").ColumnWidth = 12.576851254417766
